$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New master_device rows to append (ids 3000176-3000180)
$newRows = @(
    @{ Id = 3000176; Name = "Finger Print Scanner 32"; Mac = "80-75-40-E8-CA-24"; Serial = "BS563Q2230824"; Dspec = 165 },
    @{ Id = 3000177; Name = "IRIS Scanner 32";          Mac = "0E-1A-14-4A-6D-3A"; Serial = "BS563Q2230825"; Dspec = 327 },
    @{ Id = 3000178; Name = "Web Camera 32";             Mac = "65-13-7F-0F-F7-53"; Serial = "BS563Q2230826"; Dspec = 736 },
    @{ Id = 3000179; Name = "Document Scanner 32";       Mac = "73-C4-DE-8E-C9-8D"; Serial = "BS563Q2230827"; Dspec = 801 },
    @{ Id = 3000180; Name = "Printer 32";                Mac = "EC-74-AB-E0-0F-38"; Serial = "BS563Q2230828"; Dspec = 920 }
)

$row = 157
foreach ($r in $newRows) {
    $ws.Cells.Item($row, 1).Value = $r.Id
    $ws.Cells.Item($row, 2).Value = $r.Name
    $ws.Cells.Item($row, 3).Value = $r.Mac
    $ws.Cells.Item($row, 4).Value = $r.Serial
    $ws.Cells.Item($row, 6).Value = $r.Dspec
    $ws.Cells.Item($row, 7).Value = "eng"
    $ws.Cells.Item($row, 8).Value = $true
    $ws.Cells.Item($row, 9).Value = "superadmin"
    $ws.Cells.Item($row, 10).Value = "now()"
    $ws.Cells.Item($row, 11).Value = "now()"
    $ws.Range("H" + $row).HorizontalAlignment = -4131
    $row++
}

# Trailing empty rows (162-166) retain the left-aligned style on column H only
for ($i = 162; $i -le 166; $i++) {
    $ws.Range("H" + $i).HorizontalAlignment = -4131
}

# Update the view to reflect where the editor ended up after adding the rows
$ws.Range("E159").Select()
$excel.ActiveWindow.ScrollRow = 154
